$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the email addresses (extend each by one extra digit)
$ws.Range("C1").Value = "es2345@yahoo.com"
$ws.Range("C2").Value = "us12345@test.com"
